$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14-17 (Resolving-Mac x MuSCs target rows removed; table now 4 senders x 3 targets)
$ws.Rows("14:17").Delete()

# Update data rows 2-13 with recomputed TPM-based values
# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Icam4"
$ws.Cells.Item(2,3).Value = "Itgb2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.389221
$ws.Cells.Item(2,8).Value = 4.167663
$ws.Cells.Item(2,9).Value = 0.2910270461264192
$ws.Cells.Item(2,10).Value = 0.2910270461264192
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.1145113333333333
$ws.Cells.Item(2,14).Value = 0.343534
$ws.Cells.Item(2,15).Value = 0.001785365609625045
$ws.Cells.Item(2,16).Value = 0.001785365609625044
$ws.Cells.Item(2,17).Value = 0.1590815490046667
$ws.Cells.Item(2,18).Value = 1.431733941042
$ws.Cells.Item(2,19).Value = 0.0005195896796248704
$ws.Cells.Item(2,20).Value = 0.0005195896796248703

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Icam4"
$ws.Cells.Item(3,3).Value = "Itgb2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.389221
$ws.Cells.Item(3,8).Value = 4.167663
$ws.Cells.Item(3,9).Value = 0.2910270461264192
$ws.Cells.Item(3,10).Value = 0.2910270461264192
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.467525
$ws.Cells.Item(3,14).Value = 1.402575
$ws.Cells.Item(3,15).Value = 0.007289261528465441
$ws.Cells.Item(3,16).Value = 0.007289261528465441
$ws.Cells.Item(3,17).Value = 0.6494955480250001
$ws.Cells.Item(3,18).Value = 5.845459932225
$ws.Cells.Item(3,19).Value = 0.002121372251072244
$ws.Cells.Item(3,20).Value = 0.002121372251072244

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Icam4"
$ws.Cells.Item(4,3).Value = "Itgb2"
$ws.Cells.Item(4,4).Value = "Resolving-Mac"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.389221
$ws.Cells.Item(4,8).Value = 4.167663
$ws.Cells.Item(4,9).Value = 0.2910270461264192
$ws.Cells.Item(4,10).Value = 0.2910270461264192
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 63.556834
$ws.Cells.Item(4,14).Value = 190.670502
$ws.Cells.Item(4,15).Value = 0.9909253728619096
$ws.Cells.Item(4,16).Value = 0.9909253728619095
$ws.Cells.Item(4,17).Value = 88.29448848631401
$ws.Cells.Item(4,18).Value = 794.6503963768261
$ws.Cells.Item(4,19).Value = 0.2883860841957221
$ws.Cells.Item(4,20).Value = 0.288386084195722

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Icam4"
$ws.Cells.Item(5,3).Value = "Itgb2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.202155333333334
$ws.Cells.Item(5,8).Value = 6.606466000000001
$ws.Cells.Item(5,9).Value = 0.461328155686921
$ws.Cells.Item(5,10).Value = 0.4613281556869209
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.1145113333333333
$ws.Cells.Item(5,14).Value = 0.343534
$ws.Cells.Item(5,15).Value = 0.001785365609625045
$ws.Cells.Item(5,16).Value = 0.001785365609625044
$ws.Cells.Item(5,17).Value = 0.2521717434271112
$ws.Cells.Item(5,18).Value = 2.269545690844001
$ws.Cells.Item(5,19).Value = 0.0008236394239151772
$ws.Cells.Item(5,20).Value = 0.000823639423915177

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Icam4"
$ws.Cells.Item(6,3).Value = "Itgb2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.202155333333334
$ws.Cells.Item(6,8).Value = 6.606466000000001
$ws.Cells.Item(6,9).Value = 0.461328155686921
$ws.Cells.Item(6,10).Value = 0.4613281556869209
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.467525
$ws.Cells.Item(6,14).Value = 1.402575
$ws.Cells.Item(6,15).Value = 0.007289261528465441
$ws.Cells.Item(6,16).Value = 0.007289261528465441
$ws.Cells.Item(6,17).Value = 1.029562672216667
$ws.Cells.Item(6,18).Value = 9.266064049950002
$ws.Cells.Item(6,19).Value = 0.003362741577246588
$ws.Cells.Item(6,20).Value = 0.003362741577246588

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Icam4"
$ws.Cells.Item(7,3).Value = "Itgb2"
$ws.Cells.Item(7,4).Value = "Resolving-Mac"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.202155333333334
$ws.Cells.Item(7,8).Value = 6.606466000000001
$ws.Cells.Item(7,9).Value = 0.461328155686921
$ws.Cells.Item(7,10).Value = 0.4613281556869209
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 63.556834
$ws.Cells.Item(7,14).Value = 190.670502
$ws.Cells.Item(7,15).Value = 0.9909253728619096
$ws.Cells.Item(7,16).Value = 0.9909253728619095
$ws.Cells.Item(7,17).Value = 139.9620209628814
$ws.Cells.Item(7,18).Value = 1259.658188665932
$ws.Cells.Item(7,19).Value = 0.4571417746857593
$ws.Cells.Item(7,20).Value = 0.4571417746857592

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Icam4"
$ws.Cells.Item(8,3).Value = "Itgb2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.1506176666666667
$ws.Cells.Item(8,8).Value = 0.451853
$ws.Cells.Item(8,9).Value = 0.03155280162368235
$ws.Cells.Item(8,10).Value = 0.03155280162368235
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.1145113333333333
$ws.Cells.Item(8,14).Value = 0.343534
$ws.Cells.Item(8,15).Value = 0.001785365609625045
$ws.Cells.Item(8,16).Value = 0.001785365609625044
$ws.Cells.Item(8,17).Value = 0.01724742983355556
$ws.Cells.Item(8,18).Value = 0.155226868502
$ws.Cells.Item(8,19).Value = [double]"5.633328690624374E-05"
$ws.Cells.Item(8,20).Value = [double]"5.633328690624373E-05"

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Icam4"
$ws.Cells.Item(9,3).Value = "Itgb2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.1506176666666667
$ws.Cells.Item(9,8).Value = 0.451853
$ws.Cells.Item(9,9).Value = 0.03155280162368235
$ws.Cells.Item(9,10).Value = 0.03155280162368235
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.467525
$ws.Cells.Item(9,14).Value = 1.402575
$ws.Cells.Item(9,15).Value = 0.007289261528465441
$ws.Cells.Item(9,16).Value = 0.007289261528465441
$ws.Cells.Item(9,17).Value = 0.07041752460833334
$ws.Cells.Item(9,18).Value = 0.6337577214750001
$ws.Cells.Item(9,19).Value = 0.0002299966229908097
$ws.Cells.Item(9,20).Value = 0.0002299966229908097

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Icam4"
$ws.Cells.Item(10,3).Value = "Itgb2"
$ws.Cells.Item(10,4).Value = "Resolving-Mac"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.1506176666666667
$ws.Cells.Item(10,8).Value = 0.451853
$ws.Cells.Item(10,9).Value = 0.03155280162368235
$ws.Cells.Item(10,10).Value = 0.03155280162368235
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 63.556834
$ws.Cells.Item(10,14).Value = 190.670502
$ws.Cells.Item(10,15).Value = 0.9909253728619096
$ws.Cells.Item(10,16).Value = 0.9909253728619095
$ws.Cells.Item(10,17).Value = 9.572782037800668
$ws.Cells.Item(10,18).Value = 86.155038340206
$ws.Cells.Item(10,19).Value = 0.03126647171378531
$ws.Cells.Item(10,20).Value = 0.0312664717137853

# Row 11
$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,2).Value = "Icam4"
$ws.Cells.Item(11,3).Value = "Itgb2"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 1.031517666666667
$ws.Cells.Item(11,8).Value = 3.094553
$ws.Cells.Item(11,9).Value = 0.2160919965629775
$ws.Cells.Item(11,10).Value = 0.2160919965629775
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.1145113333333333
$ws.Cells.Item(11,14).Value = 0.343534
$ws.Cells.Item(11,15).Value = 0.001785365609625045
$ws.Cells.Item(11,16).Value = 0.001785365609625044
$ws.Cells.Item(11,17).Value = 0.1181204633668889
$ws.Cells.Item(11,18).Value = 1.063084170302
$ws.Cells.Item(11,19).Value = 0.0003858032191787534
$ws.Cells.Item(11,20).Value = 0.0003858032191787534

# Row 12
$ws.Cells.Item(12,1).Value = "Resolving-Mac"
$ws.Cells.Item(12,2).Value = "Icam4"
$ws.Cells.Item(12,3).Value = "Itgb2"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 1.031517666666667
$ws.Cells.Item(12,8).Value = 3.094553
$ws.Cells.Item(12,9).Value = 0.2160919965629775
$ws.Cells.Item(12,10).Value = 0.2160919965629775
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.467525
$ws.Cells.Item(12,14).Value = 1.402575
$ws.Cells.Item(12,15).Value = 0.007289261528465441
$ws.Cells.Item(12,16).Value = 0.007289261528465441
$ws.Cells.Item(12,17).Value = 0.4822602971083333
$ws.Cells.Item(12,18).Value = 4.340342673975
$ws.Cells.Item(12,19).Value = 0.001575151077155798
$ws.Cells.Item(12,20).Value = 0.001575151077155798

# Row 13
$ws.Cells.Item(13,1).Value = "Resolving-Mac"
$ws.Cells.Item(13,2).Value = "Icam4"
$ws.Cells.Item(13,3).Value = "Itgb2"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 1.031517666666667
$ws.Cells.Item(13,8).Value = 3.094553
$ws.Cells.Item(13,9).Value = 0.2160919965629775
$ws.Cells.Item(13,10).Value = 0.2160919965629775
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 63.556834
$ws.Cells.Item(13,14).Value = 190.670502
$ws.Cells.Item(13,15).Value = 0.9909253728619096
$ws.Cells.Item(13,16).Value = 0.9909253728619095
$ws.Cells.Item(13,17).Value = 65.55999710840067
$ws.Cells.Item(13,18).Value = 590.039973975606
$ws.Cells.Item(13,19).Value = 0.214131042266643
$ws.Cells.Item(13,20).Value = 0.214131042266643
